$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.451.52"
$ws.Range("E2").Value = "  +2.08%  "
$ws.Range("D3").Value = "1.872.81"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("D4").Value = "'1.017"
$ws.Range("E4").Value = "  +0.82%  "
$ws.Range("D5").Value = "'313.45"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").Value = "'1.015"
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("D7").Value = "'0.4792"
$ws.Range("E7").Value = "  +1.41%  "
$ws.Range("D8").Value = "'0.3762"
$ws.Range("E8").Value = "  +2.68%  "
$ws.Range("D9").Value = "'0.07377"
$ws.Range("E9").Value = "  +2.68%  "
$ws.Range("D10").Value = "'0.9409"
$ws.Range("E10").Value = "  +2.05%  "
$ws.Range("D11").Value = "'20.69"
$ws.Range("E11").Value = "  +5.69%  "
$ws.Range("D12").Value = "'0.07895"
$ws.Range("E12").Value = "  +3.67%  "
$ws.Range("D13").Value = "1.892.54"
$ws.Range("E13").Value = "  +2.98%  "
$ws.Range("D14").Value = "'5.439"
$ws.Range("E14").Value = "  +2.96%  "
$ws.Range("D15").Value = "'6.604"
$ws.Range("E15").Value = "  +3.33%  "
$ws.Range("D16").Value = "'90.83"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("D17").Value = "'1.016"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").Value = "'0.000008906"
$ws.Range("E18").Value = "  +3.11%  "
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("D20").Value = "'14.91"
$ws.Range("E20").Value = "  +2.63%  "
$ws.Range("D21").Value = "27.482.05"
$ws.Range("E21").Value = "  +2.09%  "
$ws.Range("D22").Value = "'5.148"
$ws.Range("E22").Value = "  +2.85%  "
$ws.Range("D23").Value = "'10.70"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").Value = "'1.964"
$ws.Range("E24").Value = "  +2.50%  "
$ws.Range("D25").Value = "'154.36"
$ws.Range("D26").Value = "'18.59"
$ws.Range("E26").Value = "  +2.41%  "
$ws.Range("D27").Value = "'2.015"
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("D28").Value = "'116.11"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("D29").Value = "'5.011"
$ws.Range("E29").Value = "  +3.31%  "
$ws.Range("D30").Value = "'0.08931"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("D31").Value = "'3.332"
$ws.Range("E31").Value = "  +0.85%  "
$ws.Range("D32").Value = "'1.216"
$ws.Range("E32").Value = "  +4.22%  "
$ws.Range("D33").Value = "'4.588"
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("D34").Value = "'0.7484"
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("D35").Value = "'2.694"
$ws.Range("E35").Value = "  -3.39%  "
$ws.Range("D36").Value = "'0.02068"
$ws.Range("E36").Value = "  +6.07%  "
$ws.Range("E37").Value = "  +3.38%  "
$ws.Range("D38").Value = "'0.05302"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("D39").Value = "'3.002"
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("D40").Value = "'0.5362"
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("D41").Value = "'7.077"
$ws.Range("E41").Value = "  +2.83%  "
$ws.Range("D42").Value = "'0.1527"
$ws.Range("E42").Value = "  +1.11%  "
$ws.Range("D43").Value = "'8.421"
$ws.Range("E43").Value = "  +3.10%  "
$ws.Range("D44").Value = "'10.64"
$ws.Range("E44").Value = "  +1.66%  "
$ws.Range("D45").Value = "'0.4841"
$ws.Range("E45").Value = "  +3.31%  "
$ws.Range("D46").Value = "'1.016"
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("D47").Value = "'1.663"
$ws.Range("E47").Value = "  +3.96%  "
$ws.Range("D48").Value = "'103.30"
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("D49").Value = "'67.06"
$ws.Range("E49").Value = "  +2.78%  "
$ws.Range("D50").Value = "'0.06108"
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("D51").Value = "'0.9004"
$ws.Range("E51").Value = "  +1.71%  "
